$wb = $excel.ActiveWorkbook

# Reference sheet to copy structure/styles from (UK sheet has 12 rows;
# the new market sheets need 11 rows like Belgium, so row 10 "LCD800"
# is removed after copying).
$uk = $wb.Worksheets.Item("UK")

function Add-MarketSheet {
    param([string]$SheetTitle, [string]$MarketLabel, [string]$TicketRef)

    [void]$uk.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
    $new = $wb.Worksheets.Item($wb.Worksheets.Count)
    $new.Name = $SheetTitle
    [void]$new.Rows.Item(10).Delete()
    $new.Range("B2").Value = $MarketLabel
    $new.Range("B4").Value = $TicketRef
    return $new
}

$denmark = Add-MarketSheet "Denmark" "Denmark Market" "NGC-3446/T2004/T2005"
[void]$denmark.Range("A1:XFD1048576").Select()

$sweden = Add-MarketSheet "Sweden" "Sweden Market" "NGC-3465/T2025/T2023/T2027"
[void]$sweden.Range("A1:XFD1048576").Select()

$norway = Add-MarketSheet "Norway" "Norway Market" "NGC-3464/T1919/T1920"
[void]$norway.Range("B2:B4").Select()
[void]$norway.Activate()
